# Weekly update: insert a new latest-week record for
# "Terminal La Palmera de La Serena - Poroto granado" at the top of the
# data block (row 60), pushing all subsequent weekly rows down by one,
# which grows the sheet from 114 to 115 data+header rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 60; Excel shifts rows 60:114 down to 61:115
# and extends the used range / dimension to A1:R115 automatically.
$ws.Rows.Item(60).Insert()

# Populate the newly inserted row 60 with the latest week's record.
$ws.Range("A60").Value = 8
$ws.Range("B60").Value = "Terminal La Palmera de La Serena"
$ws.Range("C60").Value = "Coquimbo"
$ws.Range("D60").Value = 44957
$ws.Range("E60").Value = 4
$ws.Range("F60").Value = 100112030
$ws.Range("G60").Value = "Poroto granado"
$ws.Range("H60").Value = "Sin especificar"
$ws.Range("I60").Value = "Primera"
$ws.Range("J60").Value = 480
$ws.Range("K60").Value = 32000
$ws.Range("L60").Value = 33000
$ws.Range("M60").Value = 32500
$ws.Range("N60").Value = "`$/malla 25 kilos"
$ws.Range("O60").Value = "Provincia del Elquí"
$ws.Range("P60").Value = 1300
$ws.Range("Q60").Value = 25
$ws.Range("R60").Value = "Hortaliza"
